$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column N (year 2023) additions ---

# Row 3 header: N3 = 2023, formatted like M3
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)
$ws.Range("N3").Value = 2023

# Row 4: update M4 value, add N4 formatted like M4
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("M4").Value = 923.8
$ws.Range("N4").Value = 583.5

# Row 5: add N5 (empty), formatted like M5
$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)

# Row 6: add N6 (empty), formatted like M6
$ws.Range("M6").Copy()
$ws.Range("N6").PasteSpecial(-4122)

# Row 7: update M7 value, add N7 formatted like M7
$ws.Range("M7").Copy()
$ws.Range("N7").PasteSpecial(-4122)
$ws.Range("M7").Value = 64.03
$ws.Range("N7").Value = 64.08

$excel.CutCopyMode = $false

# --- Row heights ---
$ws.Rows.Item(1).RowHeight = 29.25
$ws.Rows.Item(4).RowHeight = 26.25
$ws.Rows.Item(5).RowHeight = 25.5
$ws.Rows.Item(6).RowHeight = 28.5
$ws.Rows.Item(7).RowHeight = 41.25

# --- Reset selection away from the stale M14 reference ---
$ws.Range("A1").Select()
